$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 5 (pushes old rows 5-15 down to 6-16), then drop the
#     leftover blank placeholder cells Excel creates in A5:C5 so that row 5
#     ends up containing only D5:H5. ---
$ws.Rows("5").Insert()
$ws.Range("A5:C5").Clear()

# --- Column widths for the new/changed columns F, G, H. ---
$ws.Columns("F").ColumnWidth = 18
$ws.Columns("G").ColumnWidth = 22.833333333333332
$ws.Columns("H").ColumnWidth = 20.833333333333332

# --- Cell writes, issued in the same order new label strings were first
#     introduced so the shared-string table lands in the expected order. ---
$ws.Range("E5").Value = "C-Eb-Gb"
$ws.Range("D5").Value = "dim"
$ws.Range("F2").Value = "M"
$ws.Range("F3").Value = "m"
$ws.Range("F4").Value = "a"
$ws.Range("F5").Value = "d"
$ws.Range("F6").Value = "d7"
$ws.Range("F11").Value = "M7"
$ws.Range("F13").Value = "m7"
$ws.Range("F8").Value = "M4"
$ws.Range("F10").Value = "m6"
$ws.Range("F9").Value = "M6"
$ws.Range("H1").Value = "OriginalBREVELabels"
$ws.Range("F12").Value = "V"
$ws.Range("H17").Value = "(deleted-d6)"
$ws.Range("F15").Value = "hd7"
$ws.Range("G1").Value = "KPLabels(29)"
$ws.Range("F1").Value = "Style2"
$ws.Range("D1").Value = "Style1"
$ws.Range("F7").Value = "M2"
$ws.Range("F14").Value = "mM7"

$ws.Range("G2").Value = "M"
$ws.Range("H2").Value = "M"
$ws.Range("G3").Value = "m"
$ws.Range("H3").Value = "m"
$ws.Range("G5").Value = "d"
$ws.Range("H5").Value = "d"
$ws.Range("G6").Value = "d7"
$ws.Range("H6").Value = "d7"
$ws.Range("H8").Value = "M4"
$ws.Range("H9").Value = "M6"
$ws.Range("G10").Value = "m6"
$ws.Range("H10").Value = "m6"
$ws.Range("G11").Value = "M7"
$ws.Range("G12").Value = "V"
$ws.Range("H12").Value = "V"
$ws.Range("G13").Value = "m7"
$ws.Range("H13").Value = "m7"
$ws.Range("G15").Value = "hd7"
$ws.Range("F16").Value = "N"
$ws.Range("G16").Value = "N"

# --- Most G/H cells carry the (text) cell style that columns A-F already
#     have at the column level; apply it explicitly cell-by-cell here
#     (G/H has no column-level style). A handful of cells (H9, G10, H10,
#     H17) are intentionally left at the default style, matching source. ---
$styledCells = @("H1","G1","G2","H2","G3","H3","G5","H5","G6","H6","H8","G11","G12","H12","G13","H13","G15","G16")
foreach ($ref in $styledCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Active selection shown in the saved view. ---
$ws.Range("I16").Select()
